$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Text) {
    # Force the cell to keep its content as literal text (e.g. "56%")
    # instead of letting Excel auto-convert it into a percentage number.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
}

# --- Row 2 (ter. 24) ---
Set-TextValue $ws.Range("B2") "34°"
Set-TextValue $ws.Range("D2") "56%"
Set-TextValue $ws.Range("E2") "81%"

# --- Row 3 (qua. 25) ---
Set-TextValue $ws.Range("E3") "77%"

# --- Row 4 (qui. 26) ---
Set-TextValue $ws.Range("D4") "59%"

# --- Row 5 (sex. 27) ---
Set-TextValue $ws.Range("E5") "82%"

# --- Row 6 (sáb. 28) ---
Set-TextValue $ws.Range("C6") "26°"
Set-TextValue $ws.Range("D6") "65%"
Set-TextValue $ws.Range("E6") "85%"

# --- Row 7 (dom. 29) ---
Set-TextValue $ws.Range("D7") "69%"
Set-TextValue $ws.Range("E7") "85%"

# --- Row 8 (seg. 30) ---
Set-TextValue $ws.Range("B8") "35°"
Set-TextValue $ws.Range("D8") "64%"
Set-TextValue $ws.Range("E8") "86%"

# --- Row 9 (ter. 01) ---
Set-TextValue $ws.Range("D9") "73%"
Set-TextValue $ws.Range("E9") "91%"

# --- Row 10 (qua. 02) ---
Set-TextValue $ws.Range("B10") "34°"
Set-TextValue $ws.Range("D10") "69%"

# Remove the last data row (qui. 03), shrinking the table from 11 to 10 rows
$ws.Rows("11:11").Delete()
